$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 255, shifting rows 255..324 down to 256..325.
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new record.
$ws.Range("A255").Value = 3
$ws.Range("B255").Value = "Femacal de La Calera"
$ws.Range("C255").Value = "Coquimbo"
$ws.Range("D255").Value = 44642
$ws.Range("E255").Value = 5
$ws.Range("F255").Value = 100112040
$ws.Range("G255").Value = "Cilantro"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 125
$ws.Range("K255").Value = 4500
$ws.Range("L255").Value = 5000
$ws.Range("M255").Value = 4760
$ws.Range("N255").Value = "$/docena de atados (3 kilos)"
$ws.Range("O255").Value = "Provincia de Quillota"
$ws.Range("P255").Value = 1587
$ws.Range("Q255").Value = 3
$ws.Range("R255").Value = "Hortaliza"
